# [LC-943] Update documentation for Letsco OS 1.3.1
# Rename sheets GP1/GP2/BP1..BP7 to GP01/GP02/BP01..BP07 (zero-padded)
# and update the matching "KPI <name> - ..." title text in cell A1 of each sheet.

$wb = $excel.ActiveWorkbook

$map = @(
    @{ Old = "GP1"; New = "GP01" },
    @{ Old = "GP2"; New = "GP02" },
    @{ Old = "BP1"; New = "BP01" },
    @{ Old = "BP2"; New = "BP02" },
    @{ Old = "BP3"; New = "BP03" },
    @{ Old = "BP4"; New = "BP04" },
    @{ Old = "BP5"; New = "BP05" },
    @{ Old = "BP6"; New = "BP06" },
    @{ Old = "BP7"; New = "BP07" }
)

foreach ($entry in $map) {
    $oldName = $entry.Old
    $newName = $entry.New

    $ws = $wb.Worksheets.Item($oldName)

    # Update the title text held in A1 (e.g. "KPI GP1 - Global Perf 1" -> "KPI GP01 - Global Perf 1")
    $cell = $ws.Range("A1")
    $currentValue = [string]$cell.Text
    $cell.Value = $currentValue.Replace("KPI $oldName ", "KPI $newName ")

    # Rename the sheet itself last, since lookups above used the old name
    $ws.Name = $newName
}
